## DevTesting_IC.dry.xlsx edits
##
## 1) Reservoirs sheet: the frozen-pane viewport was scrolled so its visible
##    top-left cell becomes H2 (pane stays frozen), and the remembered
##    selection in the frozen (bottom-right) pane moves from A17 to L16.
## 2) CoordinatedOps sheet: the remembered selection moves from C6 to E5,
##    and a TARV (Target Annual Release Volume) value of 8,232,000 is
##    entered for the 2020 row (row 4, column E) -- replacing the blank/
##    "NaN" placeholder that was there before.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Reservoirs sheet: move the frozen pane's visible window + selection
# ---------------------------------------------------------------------
$wsRes = $wb.Worksheets.Item("Reservoirs")
$wsRes.Activate()
$win = $excel.ActiveWindow

# Re-anchor the freeze at the new scroll position (H2) so the pane's
# topLeftCell updates, then restore the selection that should be active
# in the bottom-right (frozen) pane.
$win.FreezePanes = $false
$wsRes.Range("H2").Select()
$win.FreezePanes = $true
$wsRes.Range("L16").Select()

# ---------------------------------------------------------------------
# CoordinatedOps sheet: update the TARV value and the remembered selection
# ---------------------------------------------------------------------
$wsOps = $wb.Worksheets.Item("CoordinatedOps")
$wsOps.Activate()

# Row 4 = 1/1/2020. Column E = TargetAnnualReleaseVolume_Input.
$wsOps.Range("E4").Value = 8232000

$wsOps.Range("E5").Select()
